$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the regression coefficient labels to their rounded values
$ws.Range("B2").Value = "-0.37***"
$ws.Range("B3").Value = "-3.46***"
$ws.Range("C3").Value = "-0.81***"
